# Actualización desde MV -datos-
# Adds 4 new daily rows (02-11-2021 .. 05-11-2021) to the bottom of the
# "Tasas de Spreads soberanos EMBI 2021 - Diaria" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 211; Values = @("02-11-2021", 329.8, 408, 204, 320, 82.40000000000001, 117.7, 41.3, 158.1, 501.5, 1740.7, 344.7, 308, 163, 361, 183) },
    @{ Row = 212; Values = @("03-11-2021", 324.8, 403, 198, 314, 79.59999999999999, 112.9, 40.5, 153.9, 491.4, 1737.3, 336.6, 306, 155, 357, 175) },
    @{ Row = 213; Values = @("04-11-2021", 327.5, 405, 202, 315, 83.2, 117.5, 41.5, 157.7, 484.9, 1746, 331.5, 310, 157, 358, 179) },
    @{ Row = 214; Values = @("05-11-2021", 329.1, 406, 205, 313, 84.40000000000001, 118.2, 43, 160.1, 480.4, 1725.7, 331.4, 313, 160, 362, 182) }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $values = $entry.Values

    # Column A holds the date label as plain text. Excel auto-recognizes
    # strings like "02-11-2021" as dates, so force the cell to Text format
    # before writing it, then restore the default ("Normal") style so the
    # cell matches the look of the rest of the column.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $values[0]
    $cellA.Style = "Normal"

    for ($c = 2; $c -le 16; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}
